$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark these fields as mandatory ("S" = Sim) instead of "N" (Nao),
# as part of the multi-record layout normalization/validation work.
$ws.Range("E2:E9").Value = "S"
$ws.Range("E11:E12").Value = "S"
